# Fix OLS index error, add classification and ROC curve
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-7: populate previously-empty C:G cells with computed metrics
$data = @(
    @{ Row = 2;  C = 0.575;  D = 15.999; E = 12.7;   F = 234.35;  G = 6 },
    @{ Row = 3;  C = -1.451; D = 33.237; E = 26.695; F = 232.136; G = 6 },
    @{ Row = 4;  C = -1.336; D = 38.425; E = 30.355; F = 238.328; G = 6 },
    @{ Row = 5;  C = 0.485;  D = 17.608; E = 13.866; F = 236.27;  G = 6 },
    @{ Row = 6;  C = -1.519; D = 33.697; E = 26.396; F = 232.805; G = 6 },
    @{ Row = 7;  C = -1.404; D = 38.973; E = 30.127; F = 240.285; G = 6 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = $entry.F
    $ws.Cells.Item($r, 7).Value = $entry.G
}

# Rows 8-10: clear the numeric G value back to an empty (no N_obs_test available)
foreach ($r in 8..10) {
    $ws.Cells.Item($r, 7).Value = ""
}
